$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the existing "TS9310 / 3 meses" block (old rows 9-12),
# pushing it down to rows 12-15 and leaving room for a new title row (row 11).
$ws.Rows("9:11").Insert()

# Write the new multi-line TS9311 description first so it lands at shared-string
# index 46, matching the order new strings were introduced upstream.
$ws.Range("F19").Value = "exp/TS9311/dataset_future.csv.gz #futuro`nexp/TS9311/dataset_train_final.csv.gz #para entrenar modelo final`nexp/TS9311/dataset_training.csv.gz #para hacer bo"
$ws.Range("F19").WrapText = $true
$ws.Rows(19).RowHeight = 60

# New title for the (now relocated) 3-month training experiment block.
$ws.Range("A11").Value = "EXPERMENTO CON 3 MESES DE TRAINING"

# The last row of that block (now row 15) gains a "SI" marker in column A,
# matching the other rows of the block.
$ws.Range("A15").Value = "SI"

# New section: 9-month training experiment.
$ws.Range("A18").Value = "EXPERMENTO CON 9 MESES DE TRAINING"

$ws.Range("B19").Value = "Training strategy. Separar datasets: bo, train y test"
$ws.Range("C19").Value = "Todos los campos del paso anterior"
$ws.Range("D19").Value = 9312
$ws.Range("E19").Value = "exp/FE9252/dataset.csv.gz"

$ws.Range("B20").Value = "BO"
$ws.Range("C20").Value = "Dataset del paso anterior"
$ws.Range("D20").Value = "941_HT"
$ws.Range("E20").Value = "exp/TS9310/dataset_training.csv.gz"
$ws.Range("F20").Value = "exp/HT9410/dataset_training.csv.gz"

$ws.Range("B21").Value = "Modelo final"
$ws.Range("C21").Value = "Entrenar el modelo final"
$ws.Range("D21").Value = "991_ZZ_lightgbm"
$ws.Range("E21").Value = "exp/HT9410/dataset_training.csv.gz"
$ws.Range("F21").Value = "exp/ZZ9410"

$ws.Range("B22").Value = "Modelo final"
$ws.Range("C22").Value = "Entrenar el modelo final. Cortes hasta 15000"
$ws.Range("D22").Value = "991_ZZ_lightgbm_15000"
$ws.Range("E22").Value = "exp/HT9410/dataset_training.csv.gz"
$ws.Range("F22").Value = "exp/ZZ9411"

$ws.Range("B21").Select()
